$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72 for the "quad_desc" key/value pair, which pushes
# the existing rows 72-90 down to 73-91.
$ws.Rows.Item(72).Insert()
$ws.Cells.Item(72, 1).Value = "quad_desc"
$ws.Cells.Item(72, 2).Value = "· Formed by four straight lines connected in a loop.\n· Points are plotted on a plane.\n· Angles sum up to 360°."

# Append the four new "level_intro_6_*" key/value rows at the end of the
# table (rows 92-95).
$ws.Cells.Item(92, 1).Value = "level_intro_6_0"
$ws.Cells.Item(92, 2).Value = "We'll now be looking at the next sub-category of polygons, the quadrilaterals."

$ws.Cells.Item(93, 1).Value = "level_intro_6_1"
$ws.Cells.Item(94, 1).Value = "level_intro_6_2"
$ws.Cells.Item(95, 1).Value = "level_intro_6_3"

$ws.Cells.Item(93, 2).Value = "As the name suggests, it is four-sided."
$ws.Cells.Item(94, 2).Value = "For this level, we'll be looking at three particular sub-categories that all have two pairs of opposite, parallel, equal-length sides."
$ws.Cells.Item(95, 2).Value = "Just as certain triangles have more than one matching categories, so does quadrilaterals. Keep that in mind!"

# Update the view state to match where Excel would land after this edit.
$ws.Range("B95").Select()
